# Add a "Name" row (use case Name field) above the existing Brief Description
# row, formatted like the other header cells (bold, centered) but without
# word-wrap, and nudge the selection / page setup to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new label into A3 (previously empty).
$nameCell = $ws.Range("A3")
$nameCell.Value = "Name"

# Match formatting of the other bold header cells (fontId 1 = bold Calibri)
# but centered without wrap text.
$nameCell.Font.Bold = $true
$nameCell.HorizontalAlignment = -4108   # xlCenter
$nameCell.VerticalAlignment = -4108     # xlCenter
$nameCell.WrapText = $false

# Select the newly added cell, as Excel would leave it selected after typing.
$nameCell.Select()

# Page setup, as recorded by Excel when the sheet was last saved.
$ps = $ws.PageSetup
$ps.PaperSize = 9        # xlPaperA4
$ps.Orientation = 1      # xlPortrait
